$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Range("A31").Value = "udit"
$ws.Range("B31").Value = "x"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "90009749939"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "sdlfjsldkjf"
$ws.Range("E31").Value = "divesh"
$ws.Range("F31").Value = "y"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "090909989"
$ws.Range("G31").Style = "Normal"
$ws.Range("H31").Value = "muthu street"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("I31").Value = "12.09.2019"
$ws.Range("I31").Style = "Normal"
$ws.Range("J31").NumberFormat = "@"
$ws.Range("J31").Value = "25000"
$ws.Range("J31").Style = "Normal"
$ws.Range("K31").Value = "divesh"
$ws.Range("L31").Value = "yamah"
$ws.Range("M31").Value = "y3"
$ws.Range("N31").Value = "tn-02-c-8888"
$ws.Range("O31").NumberFormat = "@"
$ws.Range("O31").Value = "20"
$ws.Range("O31").Style = "Normal"
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = "2.1"
$ws.Range("P31").Style = "Normal"
$ws.Range("Q31").NumberFormat = "@"
$ws.Range("Q31").Value = "1000"
$ws.Range("Q31").Style = "Normal"
$ws.Range("R31").NumberFormat = "@"
$ws.Range("R31").Value = "3"
$ws.Range("R31").Style = "Normal"
$ws.Range("S31").Value = "d,d,d,"
$ws.Range("T31").NumberFormat = "@"
$ws.Range("T31").Value = "123"
$ws.Range("T31").Style = "Normal"
$ws.Range("U31").Value = "51b44048-dbaf-11e9-932c-107d1a2a80c2"

# Row 32
$ws.Range("A32").Value = "divesh"
$ws.Range("B32").Value = "k"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "9884523855"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "muthu street"
$ws.Range("E32").Value = "lkj"
$ws.Range("F32").Value = "k"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "9884523855"
$ws.Range("G32").Style = "Normal"
$ws.Range("H32").Value = "neil "
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = "10.12.2019"
$ws.Range("I32").Style = "Normal"
$ws.Range("J32").NumberFormat = "@"
$ws.Range("J32").Value = "12000"
$ws.Range("J32").Style = "Normal"
$ws.Range("K32").Value = "divesh"
$ws.Range("L32").Value = "yamaha"
$ws.Range("M32").Value = "y3"
$ws.Range("N32").Value = "tn-02-cc-1299"
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = "12"
$ws.Range("O32").Style = "Normal"
$ws.Range("P32").NumberFormat = "@"
$ws.Range("P32").Value = "2"
$ws.Range("P32").Style = "Normal"
$ws.Range("Q32").NumberFormat = "@"
$ws.Range("Q32").Value = "1000"
$ws.Range("Q32").Style = "Normal"
$ws.Range("R32").NumberFormat = "@"
$ws.Range("R32").Value = "2"
$ws.Range("R32").Style = "Normal"
$ws.Range("S32").Value = "kh"
$ws.Range("T32").NumberFormat = "@"
$ws.Range("T32").Value = "123"
$ws.Range("T32").Style = "Normal"
$ws.Range("U32").Value = "bfab0d98-f55d-11e9-b250-f8da0c2bca8e"

# Row 33
$ws.Range("A33").Value = "Rahul"
$ws.Range("B33").Value = "Navaratan"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "9043926545"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "Muthukrishnan street"
$ws.Range("E33").Value = "Nilesh"
$ws.Range("F33").Value = "Navaratan"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "9884523866"
$ws.Range("G33").Style = "Normal"
$ws.Range("H33").Value = "abc colony"
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = "12.12.2019"
$ws.Range("I33").Style = "Normal"
$ws.Range("J33").NumberFormat = "@"
$ws.Range("J33").Value = "100000"
$ws.Range("J33").Style = "Normal"
$ws.Range("K33").Value = "divesh"
$ws.Range("L33").Value = "s9"
$ws.Range("M33").Value = "Yamaha"
$ws.Range("N33").Value = "tn-02-12-2234"
$ws.Range("O33").NumberFormat = "@"
$ws.Range("O33").Value = "12"
$ws.Range("O33").Style = "Normal"
$ws.Range("P33").NumberFormat = "@"
$ws.Range("P33").Value = "2.5"
$ws.Range("P33").Style = "Normal"
$ws.Range("Q33").NumberFormat = "@"
$ws.Range("Q33").Value = "1000"
$ws.Range("Q33").Style = "Normal"
$ws.Range("R33").NumberFormat = "@"
$ws.Range("R33").Value = "2"
$ws.Range("R33").Style = "Normal"
$ws.Range("S33").Value = "RC,Ration card"
$ws.Range("T33").Value = "c-12"
$ws.Range("U33").Value = "731903fa-00c8-11ea-b8ac-f8da0c2bca8e"

# Row 34
$ws.Range("A34").Value = "jatin"
$ws.Range("B34").Value = "jdasfkl"
$ws.Range("C34").Value = "klsjdf"
$ws.Range("D34").Value = "lksjdf"
$ws.Range("E34").Value = "lkjsflkj"
$ws.Range("F34").Value = "lskdjf"
$ws.Range("G34").Value = "lksdjf"
$ws.Range("H34").Value = "lksdjf"
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value = "12.08.2019"
$ws.Range("I34").Style = "Normal"
$ws.Range("J34").NumberFormat = "@"
$ws.Range("J34").Value = "30000"
$ws.Range("J34").Style = "Normal"
$ws.Range("K34").Value = "divesh"
$ws.Range("L34").Value = "yamaha"
$ws.Range("M34").Value = "y4"
$ws.Range("N34").Value = "tn-02-cc-0000"
$ws.Range("O34").NumberFormat = "@"
$ws.Range("O34").Value = "12"
$ws.Range("O34").Style = "Normal"
$ws.Range("P34").NumberFormat = "@"
$ws.Range("P34").Value = "3.1"
$ws.Range("P34").Style = "Normal"
$ws.Range("Q34").NumberFormat = "@"
$ws.Range("Q34").Value = "1000"
$ws.Range("Q34").Style = "Normal"
$ws.Range("R34").NumberFormat = "@"
$ws.Range("R34").Value = "3"
$ws.Range("R34").Style = "Normal"
$ws.Range("S34").Value = "sdkflj"
$ws.Range("T34").Value = "dsklf"
$ws.Range("U34").Value = "492bdd08-1937-11ea-b191-f8da0c2bca8e"

# Row 35
$ws.Range("A35").Value = "anuraj"
$ws.Range("B35").Value = "kldsjf"
$ws.Range("C35").Value = "lkdjsf"
$ws.Range("D35").Value = "klsjdf"
$ws.Range("E35").Value = "kldsjf"
$ws.Range("F35").Value = "dlskjf"
$ws.Range("G35").Value = "sldkfj"
$ws.Range("H35").Value = "slkdjf"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = "15.12.2019"
$ws.Range("I35").Style = "Normal"
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value = "12553"
$ws.Range("J35").Style = "Normal"
$ws.Range("K35").Value = "divesh"
$ws.Range("L35").Value = "yamaha"
$ws.Range("M35").Value = "y3"
$ws.Range("N35").Value = "tn-02-22-1234"
$ws.Range("O35").NumberFormat = "@"
$ws.Range("O35").Value = "18"
$ws.Range("O35").Style = "Normal"
$ws.Range("P35").NumberFormat = "@"
$ws.Range("P35").Value = "2.5"
$ws.Range("P35").Style = "Normal"
$ws.Range("Q35").NumberFormat = "@"
$ws.Range("Q35").Value = "1000"
$ws.Range("Q35").Style = "Normal"
$ws.Range("R35").NumberFormat = "@"
$ws.Range("R35").Value = "2"
$ws.Range("R35").Style = "Normal"
$ws.Range("S35").Value = "k"
$ws.Range("T35").Value = "hj"
$ws.Range("U35").Value = "f1ea9692-1efa-11ea-abd9-db84c14c82c2"
